$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1136.7
$ws.Range("J19").Value = 1151.7778
$ws.Range("L19").Value = 1151.7778
$ws.Range("N19").Value = -1501.7778
$ws.Range("H28").Value = 2174.5
$ws.Range("I28").Value = 1986.5
$ws.Range("K28").Value = 1986.5
$ws.Range("M28").Value = -1501.5
$ws.Range("H113").Value = 43659396
$ws.Range("I113").Value = 37039256
$ws.Range("K113").Value = 37039256
$ws.Range("M113").Value = -37036002
$ws.Range("H121").Value = 3448.5
$ws.Range("J121").Value = 3448.5
$ws.Range("L121").Value = 10345.5
$ws.Range("N121").Value = -13839.5
$ws.Range("H125").Value = 720
$ws.Range("I125").Value = 720
$ws.Range("K125").Value = 6480
$ws.Range("M125").Value = -4020
$ws.Range("H132").Value = 2691.0435
$ws.Range("I132").Value = 2585.8667
$ws.Range("K132").Value = 7757.6001
$ws.Range("M132").Value = -5227.6001
$ws.Range("H135").Value = 294756.7
$ws.Range("I135").Value = 345350.25
$ws.Range("K135").Value = 3108152.25
$ws.Range("M135").Value = -3105617.25
$ws.Range("H138").Value = 3812.6
$ws.Range("I138").Value = 896.5714
$ws.Range("K138").Value = 2689.7142
$ws.Range("M138").Value = 2450.2858

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6627
$ws.Range("J2").Value = 6627
$ws.Range("L2").Value = 6627
$ws.Range("N2").Value = -6853
$ws.Range("H32").Value = 1925086.1
$ws.Range("I32").Value = 1986034.1
$ws.Range("K32").Value = 1986034.1
$ws.Range("M32").Value = -1985747.1
$ws.Range("H45").Value = 3459.0833
$ws.Range("I45").Value = 1673.8889
$ws.Range("K45").Value = 1673.8889
$ws.Range("M45").Value = -1296.8889
$ws.Range("H61").Value = 2837.5557
$ws.Range("I61").Value = 2208.3225
$ws.Range("J61").Value = 6738.8
$ws.Range("K61").Value = 2208.3225
$ws.Range("L61").Value = 6738.8
$ws.Range("M61").Value = -1996.3225
$ws.Range("N61").Value = -7162.8
$ws.Range("H110").Value = 1210.0714
$ws.Range("I110").Value = 1210.0714
$ws.Range("K110").Value = 1210.0714
$ws.Range("M110").Value = 834.9286
$ws.Range("H116").Value = 6627
$ws.Range("J116").Value = 6627
$ws.Range("L116").Value = 6627
$ws.Range("N116").Value = -11215
$ws.Range("H122").Value = 100852.25
$ws.Range("I122").Value = 133831.67
$ws.Range("J122").Value = 1914
$ws.Range("K122").Value = 401495.01
$ws.Range("L122").Value = 5742
$ws.Range("M122").Value = -399045.01
$ws.Range("N122").Value = -10642
$ws.Range("H132").Value = 9713.808000000001
$ws.Range("I132").Value = 9023.77
$ws.Range("K132").Value = 27071.31
$ws.Range("M132").Value = -24541.31
$ws.Range("H136").Value = 2837.5557
$ws.Range("I136").Value = 2208.3225
$ws.Range("J136").Value = 6738.8
$ws.Range("K136").Value = 6624.967500000001
$ws.Range("L136").Value = 20216.4
$ws.Range("M136").Value = -4074.967500000001
$ws.Range("N136").Value = -25316.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6627
$ws.Range("J3").Value = 6627
$ws.Range("L3").Value = 6627
$ws.Range("N3").Value = -6855
$ws.Range("H94").Value = 2735.3809
$ws.Range("I94").Value = 884.5
$ws.Range("K94").Value = 884.5
$ws.Range("M94").Value = -433.5
$ws.Range("H134").Value = 8397.482
$ws.Range("I134").Value = 4637.636
$ws.Range("K134").Value = 13912.908
$ws.Range("M134").Value = -11377.908

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5958.8
$ws.Range("J16").Value = 7128.3076
$ws.Range("L16").Value = 7128.3076
$ws.Range("N16").Value = -7702.3076
$ws.Range("H31").Value = 8701.479499999999
$ws.Range("I31").Value = 3564
$ws.Range("J31").Value = 11036.697
$ws.Range("K31").Value = 3564
$ws.Range("L31").Value = 11036.697
$ws.Range("M31").Value = -3269
$ws.Range("N31").Value = -11626.697
$ws.Range("H34").Value = 8701.479499999999
$ws.Range("I34").Value = 3564
$ws.Range("J34").Value = 11036.697
$ws.Range("K34").Value = 3564
$ws.Range("L34").Value = 11036.697
$ws.Range("M34").Value = -3362
$ws.Range("N34").Value = -11440.697
$ws.Range("H56").Value = 64996.668
$ws.Range("J56").Value = 64996.668
$ws.Range("L56").Value = 64996.668
$ws.Range("N56").Value = -66686.66800000001
$ws.Range("H113").Value = 5958.8
$ws.Range("J113").Value = 7128.3076
$ws.Range("L113").Value = 7128.3076
$ws.Range("N113").Value = -11468.3076
$ws.Range("H134").Value = 8332.633
$ws.Range("I134").Value = 8318.652
$ws.Range("K134").Value = 24955.956
$ws.Range("M134").Value = -22420.956

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2002543.5
$ws.Range("I5").Value = 3334591.5
$ws.Range("J5").Value = 4471.25
$ws.Range("K5").Value = 10003774.5
$ws.Range("L5").Value = 13413.75
$ws.Range("M5").Value = -10003662.5
$ws.Range("N5").Value = -13637.75
$ws.Range("H14").Value = 13889532
$ws.Range("I14").Value = 13889532
$ws.Range("K14").Value = 41668596
$ws.Range("M14").Value = -41668423
$ws.Range("H17").Value = 2469.8572
$ws.Range("I17").Value = 1296
$ws.Range("J17").Value = 3350.25
$ws.Range("K17").Value = 3888
$ws.Range("L17").Value = 10050.75
$ws.Range("M17").Value = -3719
$ws.Range("N17").Value = -10388.75
$ws.Range("H56").Value = 6475.75
$ws.Range("I56").Value = 6475.75
$ws.Range("K56").Value = 6475.75
$ws.Range("M56").Value = -5945.75
$ws.Range("H102").Value = 9000
$ws.Range("J102").Value = 9000
$ws.Range("L102").Value = 27000
$ws.Range("N102").Value = -31868
$ws.Range("H135").Value = 2002543.5
$ws.Range("I135").Value = 3334591.5
$ws.Range("J135").Value = 4471.25
$ws.Range("K135").Value = 30011323.5
$ws.Range("L135").Value = 40241.25
$ws.Range("M135").Value = -30008788.5
$ws.Range("N135").Value = -45311.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3368.2
$ws.Range("I80").Value = 3027.5715
$ws.Range("J80").Value = 4163
$ws.Range("K80").Value = 3027.5715
$ws.Range("L80").Value = 4163
$ws.Range("M80").Value = -2029.5715
$ws.Range("N80").Value = -6159
$ws.Range("H83").Value = 3368.2
$ws.Range("I83").Value = 3027.5715
$ws.Range("J83").Value = 4163
$ws.Range("K83").Value = 15137.8575
$ws.Range("L83").Value = 20815
$ws.Range("M83").Value = -10145.8575
$ws.Range("N83").Value = -30799
$ws.Range("H97").Value = 2290.0667
$ws.Range("I97").Value = 1881.1578
$ws.Range("K97").Value = 1881.1578
$ws.Range("M97").Value = -1385.1578
$ws.Range("H133").Value = 79280
$ws.Range("J133").Value = 79280
$ws.Range("L133").Value = 79280
$ws.Range("N133").Value = -89400

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1282424.1
$ws.Range("I82").Value = 1410426.6
$ws.Range("J82").Value = 2400
$ws.Range("K82").Value = 1410426.6
$ws.Range("L82").Value = 2400
$ws.Range("M82").Value = -1410065.6
$ws.Range("N82").Value = -3122
$ws.Range("H85").Value = 1282424.1
$ws.Range("I85").Value = 1410426.6
$ws.Range("J85").Value = 2400
$ws.Range("K85").Value = 1410426.6
$ws.Range("L85").Value = 2400
$ws.Range("M85").Value = -1409178.6
$ws.Range("N85").Value = -4896
$ws.Range("H93").Value = 7764.1816
$ws.Range("I93").Value = 7267.3335
$ws.Range("K93").Value = 7267.3335
$ws.Range("M93").Value = -6019.3335
$ws.Range("H122").Value = 3666.5813
$ws.Range("J122").Value = 6213
$ws.Range("L122").Value = 18639
$ws.Range("N122").Value = -23539
$ws.Range("H136").Value = 8398.323
$ws.Range("I136").Value = 2149.6667
$ws.Range("J136").Value = 10647.84
$ws.Range("K136").Value = 6449.000100000001
$ws.Range("L136").Value = 31943.52
$ws.Range("M136").Value = -3899.000100000001
$ws.Range("N136").Value = -37043.52

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 451356
$ws.Range("I122").Value = 2003602
$ws.Range("J122").Value = 7857.143
$ws.Range("K122").Value = 6010806
$ws.Range("L122").Value = 23571.429
$ws.Range("M122").Value = -6008356
$ws.Range("N122").Value = -28471.429
$ws.Range("H132").Value = 13519352
$ws.Range("I132").Value = 20839436
$ws.Range("K132").Value = 62518308
$ws.Range("M132").Value = -62515778
$ws.Range("H136").Value = 23837010
$ws.Range("J136").Value = 45717.68
$ws.Range("L136").Value = 137153.04
$ws.Range("N136").Value = -142253.04
